# Update cryptos list: price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.210.49"
$ws.Range("E2").Value = "  +3.58%  "

$ws.Range("D3").Value = "3.114.23"
$ws.Range("E3").Value = "  +1.54%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.01%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +1.52%  "

$ws.Range("E9").Value = "  +3.08%  "

$ws.Range("E10").Value = "  +1.55%  "

$ws.Range("E11").Value = "  +3.78%  "

$ws.Range("D12").Value = "3.646.18"
$ws.Range("E12").Value = "  +1.43%  "

$ws.Range("E13").Value = "  +1.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.52%  "

$ws.Range("E15").Value = "  +2.21%  "

$ws.Range("D16").Value = "59.138.31"
$ws.Range("E16").Value = "  +3.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.82%  "

$ws.Range("D18").Value = "3.116.80"
$ws.Range("E18").Value = "  +1.31%  "

$ws.Range("E19").Value = "  +0.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "341.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.81%  "

$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("E23").Value = "  +2.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.34%  "

$ws.Range("E25").Value = "  +2.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("D27").Value = "0.0₃0926"
$ws.Range("E27").Value = "  -1.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.85%  "

$ws.Range("E29").Value = "  +4.62%  "

$ws.Range("E30").Value = "  +2.47%  "

$ws.Range("E31").Value = "  +4.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "155.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.54%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "27.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.66%  "

$ws.Range("E37").Value = "  +6.47%  "

$ws.Range("E38").Value = "  +2.64%  "

$ws.Range("E39").Value = "  +3.52%  "

$ws.Range("D40").Value = "3.153.75"
$ws.Range("E40").Value = "  +1.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("E43").Value = "  -0.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.45%  "

$ws.Range("D45").Value = "2.290.90"
$ws.Range("E45").Value = "  +2.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0259"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.964"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.763"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "261.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.45%  "
